$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "release/8.0.5"
$ws.Range("B8").Value = "X"
$ws.Range("C8").Value = "X"
$ws.Range("D8").Value = "X"
$ws.Range("E8").Value = "X"
